$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.668.06"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.445.47"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.03%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'572.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.31%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'159.17"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.81%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'3.443.04"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.17%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.581"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.77%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "  +0.15%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.122"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.61%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "  -0.99%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'4.036.83"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.28%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "  -0.11%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'27.53"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.13%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "  -8.70%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'64.663.75"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.16%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.443.45"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.64%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'6.18"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.60%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'13.78"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.39%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'380.53"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.61%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'7.98"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.24%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.21%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'72.43"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.37%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.531"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.98%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'9.91"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.41%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "  +0.67%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'0.988"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.05%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'6.08"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.19%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "  -3.71%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'2.02"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.66%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'23.24"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.02%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'7.02"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.94%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.58"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.83%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'161.23"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'1.88"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.79%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'2.880.94"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.51%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.0746"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.27%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'26.36"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.14%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "  +2.90%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'4.53"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.09%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "  +0.46%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'6.50"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.78%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'26.01"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.30%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "  -2.13%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'2.43"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +13.60%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'322.55"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.91%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "  -1.94%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "  -2.32%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.845"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.70%  "
$ws.Range("E51").Style = "Normal"
